$d = $word.ActiveDocument

# --- Change 1: "status" field documentation - number -> string, new enum values ---
$old1 = "status (number): Status of the light. (0 " + [char]0x2013 + " OFF, 1 " + [char]0x2013 + " ON, 2 " + [char]0x2013 + " NORMAL)"
$new1 = "status (string): Status of the light. (" + [char]0x201C + "Auto" + [char]0x201D + ", " + [char]0x201C + "Force_On" + [char]0x201D + ", " + [char]0x201C + "Force_Off" + [char]0x201D + ")"
$r1 = $d.Content
$found1 = $r1.Find.Execute($old1, $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Change 1 found: $found1"
if ($found1) {
    $r1.Text = $new1
}

# --- Change 2: "msSens" field documentation - string -> number, add range note ---
$old2 = "msSens (string): Motion sensor sensitivity."
$new2 = "msSens (number): Motion sensor sensitivity. (1-5)"
$r2 = $d.Content
$found2 = $r2.Find.Execute($old2, $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Change 2 found: $found2"
if ($found2) {
    $r2.Text = $new2
}

# --- Change 3: JSON example payload - status value 1 -> "Force_On" ---
$old3 = '          "status": 1,'
$new3 = '          "status": ' + [char]0x201C + "Force_On" + [char]0x201D + ","
$r3 = $d.Content
$found3 = $r3.Find.Execute($old3, $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Change 3 found: $found3"
if ($found3) {
    $r3.Text = $new3
}
